$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column R (18th column) extends the year series in row 4, and the
# corresponding data rows 5 and 6, mirroring the formatting of column Q
# (xlPasteFormats = -4122, so the numeric value we set afterwards isn't
# clobbered by the copy).

$xlPasteFormats = -4122

$ws.Range("Q4").Copy() | Out-Null
$ws.Range("R4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("R4").Value = 2022

$ws.Range("Q5").Copy() | Out-Null
$ws.Range("R5").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("R5").Value = 8.6821914120339212

$ws.Range("Q6").Copy() | Out-Null
$ws.Range("R6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("R6").Value = 12.221423436376707

$excel.CutCopyMode = $false

# Move the active selection to match the post-edit workbook state.
$ws.Range("S4").Select() | Out-Null
